$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 173625
$ws.Range("C4").Value = 163873
$ws.Range("C5").Value = 9752
$ws.Range("C6").Value = 775
$ws.Range("C7").Value = 5.62

$wb.Save()
